# Add Jing Chen's wk8 Diaries
# Fills rows 19-22 of the diary worksheet with four new entries, reusing the
# existing "filled diary row" formatting (rows 10-18) as the template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 19 -----------------------------------------------------------
$ws.Range("A16:G16").Copy()
$ws.Range("A19:G19").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("G16").Copy()
$ws.Range("D19").PasteSpecial(-4122)       # Goal cell uses the narrower italic style here

$ws.Range("A19").Value = 43888
$ws.Range("B19").Value = "17:00-20:00"
$ws.Range("D19").Value = "Another three Key expert practices & design patterns"
$ws.Range("E19").Value = "Gained knowledges about the three Key expert practices & some design patterns"
$ws.Range("F19").Value = "Useful skills and knowledge. Design patterns are popular currently and they are indeed important. I will learn more about them and try to use them in practice to improve my own project."
$ws.Range("G19").Value = "Great guest talking. Good to know about someone who are excited to develope software for doing scientific research."

# --- Row 20 -----------------------------------------------------------
$ws.Range("A16:G16").Copy()
$ws.Range("A20:G20").PasteSpecial(-4122)

$ws.Range("A20").Value = 43889
$ws.Range("B20").Value = "20:00-22:00"
$ws.Range("D20").Value = "Look for five design patterns and open issues we can fix in our program. Record them."
$ws.Range("E20").Value = "Found five design patterns and several open issues we can fix in our program."
$ws.Range("F20").Value = "Some issues may be easier to be fixed than I expected"
$ws.Range("G20").Value = "Not bad"

# --- Row 21 -----------------------------------------------------------
$ws.Range("A16:G16").Copy()
$ws.Range("A21:G21").PasteSpecial(-4122)
$ws.Range("G16").Copy()
$ws.Range("D21").PasteSpecial(-4122)       # Goal cell uses the narrower italic style here

$ws.Range("A21").Value = 43891
$ws.Range("B21").Value = "9:00-9:30"
$ws.Range("C21").Value = "Guowei Li, Dongxin Xiang"
$ws.Range("D21").Value = "Decide what design patterns to use for the assignment and open issues to fix in our program together. Discuss other related problems."
$ws.Range("E21").Value = "Decided five patterns to be used and two open issues to be fixed."
$ws.Range("F21").Value = "Design patterns in different levels have different usage, which makes them handy and popular."
$ws.Range("G21").Value = "Discussion is great"

# --- Row 22 -----------------------------------------------------------
$ws.Range("A16:G16").Copy()
$ws.Range("A22:G22").PasteSpecial(-4122)

$ws.Range("A22").Value = 43892
$ws.Range("B22").Value = "20:00-21:00"
$ws.Range("C22").Value = "Guowei Li, Dongxin Xiang"
$ws.Range("D22").Value = "Review what we have done for the assignment and discuss what we should improve"
$ws.Range("E22").Value = "We have completed the draft and filed all the problems we have not think through. We will go to find more information by ourselves, and then we will discuss again. And we plan to consult Kaj on Wednesday."
$ws.Range("F22").Value = "Made more things clear."
$ws.Range("G22").Value = "Not bad"

# --- View bookkeeping (best-effort, cosmetic) --------------------------
$ws.Range("E28").Select()
